# The "Tim Operasi" members table (the one immediately following the
# ${tim_operasi_section} placeholder) is missing an explicit <w:tblBorders>
# block, unlike its sibling tables in this template. Add one that turns
# every edge (top/left/bottom/right/insideH/insideV) off, matching the
# <w:tblBorders> already present on the other tables in the document.

$d = $word.ActiveDocument
$tbl = $d.Tables(3)

# Disable/clear the borders first so Word seeds sz/space/color defaults
# for every edge (including the inside horizontal/vertical ones, which
# aren't reachable through Borders.DistanceFrom*).
$tbl.Borders.Enable = $false

# wdBorderTop=-1, wdBorderLeft=-2, wdBorderBottom=-3, wdBorderRight=-4,
# wdBorderHorizontal=-5, wdBorderVertical=-6
for ($i = -1; $i -ge -6; $i--) {
    $border = $tbl.Borders($i)
    $border.LineWidth = 0            # wdLineWidth0 -> w:sz="0"
    $border.LineStyle = 0            # wdLineStyleNone -> w:val="none"
}
